# Add six new weekly-scoreboard rows (306-311) to Sheet1, mirroring the
# format of the existing data rows, then move the active selection past
# the newly-added data (matches the author's post-edit cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the last existing data row (305) as a formatting template so the new
# rows pick up the same number formats (notably the date style on column B)
# without creating any new style entries.
$ws.Range("A305:M305").Copy()
$ws.Range("A306:M311").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    # Participant,  Date,   Workout Type, Duration, Distance, Elevation, Z1, Z2, Z3, Z4, Z5, Workout Level,    Week
    @("Jeremiah", 45503, "Workout", 45,  0,     0,   12,  27, 5,  0,  0,  "Agile Antelope", 8),
    @("Eric",     45503, "Run",     73,  7.27,  230, 0,   2,  5,  28, 36, "Agile Antelope", 8),
    @("Steven",   45503, "Walk",    111, 5.3,   348, 109, 2,  0,  0,  0,  "Brave Leopard",  8),
    @("Matt",     45503, "Ride",    45,  13.97, 0,   7,   38, 0,  0,  0,  "Agile Antelope", 8),
    @("Jeremiah", 45504, "Workout", 54,  0,     0,   48,  7,  0,  0,  0,  "Agile Antelope", 8),
    @("Eric",     45504, "Workout", 85,  0,     0,   51,  32, 2,  0,  0,  "Agile Antelope", 8)
)

$r = 306
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $r = $r + 1
}

# Move the frozen-pane view / selection down to just past the new data,
# matching the author's cursor position after entering the rows.
$ws.Range("A300").Select() | Out-Null
$ws.Range("A312").Select() | Out-Null
